# Adds the next diagonal of "matched errors" values to rows 16-24,
# mirroring the preprocessing step that extends the ifo GDP component
# error series by one more forecast horizon column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 16; Col = "J"; Value = 0.2150495036779461 },
    @{ Row = 17; Col = "I"; Value = 0.24 },
    @{ Row = 18; Col = "H"; Value = 0.3087982760018804 },
    @{ Row = 19; Col = "G"; Value = 0.32 },
    @{ Row = 20; Col = "F"; Value = 0.4476495795507702 },
    @{ Row = 21; Col = "E"; Value = 0.1088966743764388 },
    @{ Row = 22; Col = "D"; Value = 0.1461563307127136 },
    @{ Row = 23; Col = "C"; Value = 0.09547648014918764 },
    @{ Row = 24; Col = "B"; Value = 0.0959495356205764 }
)

foreach ($u in $updates) {
    $addr = "$($u.Col)$($u.Row)"
    $ws.Range($addr).Value = $u.Value
}
